# Adds a new "2022" column (column S) to the maternal-mortality-rate table,
# mirroring the formatting of the existing 2021 column (R) and filling in
# the new year's data, then moves the active selection down to S16 (just
# below the table), matching the post-edit state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (and provisional values) of the 2021 column (R3:R14,
# the header-separator row through the last data row) into the new column
# S so the new column inherits identical borders/fonts/number formats.
$ws.Range("R3:R14").Copy($ws.Range("S3:S14"))

# Header year for the new column.
$ws.Range("S4").Value = 2022

# Maternal mortality rate values for 2022, one per region row.
$ws.Range("S5").Value = 27.292394741221504
$ws.Range("S6").Value = 36.613942589338023
$ws.Range("S7").Value = 14.18691257315127
$ws.Range("S8").Value = 55.377118174770182
$ws.Range("S9").Value = 42.247570764681029
$ws.Range("S10").Value = 30.18817294468856
$ws.Range("S11").Value = 97.03085581214826
$ws.Range("S12").Value = 25.2
$ws.Range("S13").Value = 21.849963583394029

# Last region has no data for 2022 - shown as a dash, like other
# missing-data cells elsewhere in the table.
$ws.Range("S14").Value = "-"

# Move the selection to below the (now one column wider) table.
$ws.Range("S16").Select() | Out-Null
